$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new translation columns
$ws.Range("J1").Value = "Onkelos"
$ws.Range("K1").Value = "Jonathan"

# Match the header formatting (bold, border, centered) used by the other headers
$ws.Range("I1").Copy()
$ws.Range("J1:K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Exodus row (row 2)
$ws.Range("J2").Value = "“Go, gather the elders of Yisrael, and say to them, ‘Adonoy, the God of your fathers appeared [<b>became revealed</b>] to me—the God of Avraham, Yitzchok and Yaakov—saying, “I have indeed been mindful of you, regarding that which is being done to you in Egypt."
$ws.Range("K2").Value = "Go, and assemble the elders of Israel, and say to them, The Lord God of your fathers hath appeared unto me, the God of Abraham, Izhak, and Jakob, saying, Remembering, I have remembered you, and the injury that is done you in Mizraim;"

# Deuteronomy row (row 3)
$ws.Range("J3").Value = "And there has not ever arisen a prophet within Yisroel like Moshe, whom Adonoy knew [<b>appeared to</b>] face-to-face."
$ws.Range("K3").Value = "But no prophet hath again risen in Israel like unto Mosheh, because the Word of the Lord had known him to speak with him word for word,"
